# Update LR-pair statistics for rows 2-13 (columns E:T) per updated Natmi analysis
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @{ "E"=3; "F"=1; "G"=3.081251; "H"=9.243753; "I"=0.6013642694204734; "J"=0.6013642694204734; "K"=3; "L"=1; "M"=31.22896466666667; "N"=93.686894; "O"=0.2877106972998646; "P"=0.2877106972998646; "Q"=96.22427860813133; "R"=866.0185074731819; "S"=0.173018933286188; "T"=0.173018933286188 }
    3 = @{ "E"=3; "F"=1; "G"=3.081251; "H"=9.243753; "I"=0.6013642694204734; "J"=0.6013642694204734; "K"=3; "L"=1; "M"=40.44578266666667; "N"=121.337348; "O"=0.3726247238124506; "P"=0.3726247238124505; "Q"=124.6236082874493; "R"=1121.612474587044; "S"=0.22408319480348; "T"=0.22408319480348 }
    4 = @{ "E"=3; "F"=1; "G"=3.081251; "H"=9.243753; "I"=0.6013642694204734; "J"=0.6013642694204734; "K"=3; "L"=1; "M"=25.36964133333333; "N"=76.108924; "O"=0.2337290805561598; "P"=0.2337290805561598; "Q"=78.17023272797466; "R"=703.532094551772; "S"=0.140556317770974; "T"=0.140556317770974 }
    5 = @{ "E"=3; "F"=1; "G"=3.081251; "H"=9.243753; "I"=0.6013642694204734; "J"=0.6013642694204734; "K"=3; "L"=1; "M"=11.49855033333333; "N"=34.495651; "O"=0.1059354983315251; "P"=0.1059354983315251; "Q"=35.42991971313367; "R"=318.869277418203; "S"=0.0637058235598314; "T"=0.0637058235598314 }
    6 = @{ "E"=3; "F"=1; "G"=0.6679959999999999; "H"=2.003988; "I"=0.1303720230892577; "J"=0.1303720230892577; "K"=3; "L"=1; "M"=31.22896466666667; "N"=93.686894; "O"=0.2877106972998646; "P"=0.2877106972998646; "Q"=20.86082348147466; "R"=187.747411333272; "S"=0.03750942567140439; "T"=0.03750942567140439 }
    7 = @{ "E"=3; "F"=1; "G"=0.6679959999999999; "H"=2.003988; "I"=0.1303720230892577; "J"=0.1303720230892577; "K"=3; "L"=1; "M"=40.44578266666667; "N"=121.337348; "O"=0.3726247238124506; "P"=0.3726247238124505; "Q"=27.01762103820266; "R"=243.1585893438239; "S"=0.0485798390965051; "T"=0.0485798390965051 }
    8 = @{ "E"=3; "F"=1; "G"=0.6679959999999999; "H"=2.003988; "I"=0.1303720230892577; "J"=0.1303720230892577; "K"=3; "L"=1; "M"=25.36964133333333; "N"=76.108924; "O"=0.2337290805561598; "P"=0.2337290805561598; "Q"=16.94681893210133; "R"=152.521370388912; "S"=0.03047173308689865; "T"=0.03047173308689865 }
    9 = @{ "E"=3; "F"=1; "G"=0.6679959999999999; "H"=2.003988; "I"=0.1303720230892577; "J"=0.1303720230892577; "K"=3; "L"=1; "M"=11.49855033333333; "N"=34.495651; "O"=0.1059354983315251; "P"=0.1059354983315251; "Q"=7.680985628465333; "R"=69.128870656188; "S"=0.01381102523444962; "T"=0.01381102523444962 }
    10 = @{ "E"=3; "F"=1; "G"=1.374521; "H"=4.123563; "I"=0.2682637074902688; "J"=0.2682637074902689; "K"=3; "L"=1; "M"=31.22896466666667; "N"=93.686894; "O"=0.2877106972998646; "P"=0.2877106972998646; "Q"=42.92486774259133; "R"=386.3238096833219; "S"=0.07718233834227216; "T"=0.07718233834227217 }
    11 = @{ "E"=3; "F"=1; "G"=1.374521; "H"=4.123563; "I"=0.2682637074902688; "J"=0.2682637074902689; "K"=3; "L"=1; "M"=40.44578266666667; "N"=121.337348; "O"=0.3726247238124506; "P"=0.3726247238124505; "Q"=55.59357763676933; "R"=500.3421987309239; "S"=0.09996168991246546; "T"=0.09996168991246547 }
    12 = @{ "E"=3; "F"=1; "G"=1.374521; "H"=4.123563; "I"=0.2682637074902688; "J"=0.2682637074902689; "K"=3; "L"=1; "M"=25.36964133333333; "N"=76.108924; "O"=0.2337290805561598; "P"=0.2337290805561598; "Q"=34.87110477513466; "R"=313.839942976212; "S"=0.06270102969828713; "T"=0.06270102969828714 }
    13 = @{ "E"=3; "F"=1; "G"=1.374521; "H"=4.123563; "I"=0.2682637074902688; "J"=0.2682637074902689; "K"=3; "L"=1; "M"=11.49855033333333; "N"=34.495651; "O"=0.1059354983315251; "P"=0.1059354983315251; "Q"=15.80499890272367; "R"=142.244990124513; "S"=0.02841864953724412; "T"=0.02841864953724413 }
}

foreach ($row in $rowData.Keys) {
    foreach ($col in $rowData[$row].Keys) {
        $ws.Range("$col$row").Value = $rowData[$row][$col]
    }
}
